$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 2 ("2000年" row); rows 3 and 4 shift up to become rows 2 and 3
$ws.Rows.Item(2).Delete()
